# Add the two new "Equipment Tray" print requests to the July 2018 log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("July 2018")

# Row 15: Equipment Tray Positive
$ws.Range("A15").Value = "25-07-2018"
$ws.Range("B15").Value = "25-07-2018"
$ws.Range("C15").Value = "Equipment Tray Positive"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = "Polylite"
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 20
$ws.Range("H15").Value = 0.2
$ws.Range("I15").Value = "NA"

# Row 16: Equipment Tray Clip
$ws.Range("A16").Value = "25-07-2018"
$ws.Range("C16").Value = "Equipment Tray Clip"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = "Polylite"
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = 0.2
$ws.Range("I16").Value = "NA"

# Match the author's final cell selection recorded in the sheet view.
$ws.Range("C17").Select() | Out-Null
